$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.341.27'
$ws.Range('E2').Value = '  -3.75%  '
$ws.Range('D3').Value = '1.666.15'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.56'
$ws.Range('E5').Value = '  -2.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5164'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.008'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.06448'
$ws.Range('E8').Value = '  -2.30%  '
$ws.Range('E9').Value = '  -3.75%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.97'
$ws.Range('E10').Value = '  -4.31%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07657'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.327'
$ws.Range('E12').Value = '  -5.38%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.663.43'
$ws.Range('E13').Value = '  -2.62%  '
$ws.Range('D14').Value = '1.895.40'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5536'
$ws.Range('E15').Value = '  -3.70%  '
$ws.Range('D16').Value = '0.0₅8047'
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.58'
$ws.Range('E17').Value = '  -4.76%  '
$ws.Range('D18').Value = '26.378.40'
$ws.Range('E18').Value = '  -3.57%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '210.55'
$ws.Range('E20').Value = '  -2.37%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.400'
$ws.Range('E21').Value = '  -5.90%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.12'
$ws.Range('E22').Value = '  -3.35%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.895'
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.008'
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.72'
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.746'
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('E27').Value = '  -4.07%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.996'
$ws.Range('E28').Value = '  -3.98%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.78'
$ws.Range('E29').Value = '  -3.36%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05272'
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.263'
$ws.Range('E31').Value = '  -2.46%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.374'
$ws.Range('E32').Value = '  -3.64%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.218'
$ws.Range('E33').Value = '  -6.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.568'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.756'
$ws.Range('E35').Value = '  -4.25%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.376'
$ws.Range('E36').Value = '  -1.68%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9280'
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5731'
$ws.Range('E38').Value = '  -2.36%  '
$ws.Range('D39').Value = '1.149.49'
$ws.Range('E39').Value = '  +10.06%  '
$ws.Range('E40').Value = '  -1.98%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.008'
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8440'
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.660'
$ws.Range('E43').Value = '  -3.54%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '99.98'
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('D45').Value = '1.805.29'
$ws.Range('E45').Value = '  -2.53%  '
$ws.Range('D46').Value = '0.0₈113'
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4502'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '56.08'
$ws.Range('E48').Value = '  -3.47%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.008'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.947'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05106'
